{"js": "// Update the worksheet date and the 25 two-digit multiplication prompts.\n// Each entry is [oldText, newText]; every oldText is unique in the document\n// (verified against the canonical before/after OOXML), so a scoped\n// search-and-replace on each pair reproduces the diff exactly.\nconst pairs = [\n  [\"2023-12-20 Wednesday\", \"2023-12-21 Thursday\"],\n  [\"43\u00d716=\", \"41\u00d755=\"],\n  [\"17\u00d724=\", \"45\u00d719=\"],\n  [\"41\u00d752=\", \"63\u00d762=\"],\n  [\"26\u00d768=\", \"64\u00d724=\"],\n  [\"83\u00d784=\", \"35\u00d793=\"],\n  [\"89\u00d747=\", \"61\u00d724=\"],\n  [\"46\u00d718=\", \"32\u00d735=\"],\n  [\"69\u00d798=\", \"66\u00d776=\"],\n  [\"16\u00d719=\", \"41\u00d738=\"],\n  [\"35\u00d769=\", \"29\u00d796=\"],\n  [\"82\u00d723=\", \"16\u00d776=\"],\n  [\"58\u00d789=\", \"46\u00d778=\"],\n  [\"70\u00d728=\", \"64\u00d714=\"],\n  [\"67\u00d755=\", \"54\u00d775=\"],\n  [\"27\u00d730=\", \"77\u00d784=\"],\n  [\"94\u00d736=\", \"92\u00d797=\"],\n  [\"89\u00d711=\", \"98\u00d718=\"],\n  [\"88\u00d773=\", \"55\u00d712=\"],\n  [\"93\u00d739=\", \"39\u00d776=\"],\n  [\"73\u00d723=\", \"66\u00d729=\"],\n  [\"72\u00d724=\", \"55\u00d750=\"],\n  [\"81\u00d779=\", \"93\u00d783=\"],\n  [\"46\u00d717=\", \"94\u00d791=\"],\n  [\"34\u00d728=\", \"30\u00d741=\"],\n  [\"56\u00d795=\", \"41\u00d790=\"],\n];\n\nfor (const [oldText, newText] of pairs) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 two-digit multiplication prompts.\n# Each entry is (oldText, newText); every oldText is unique in the document\n# (verified against the canonical before/after OOXML), so a scoped\n# Find/Replace on each pair reproduces the diff exactly.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2023-12-20 Wednesday\", \"2023-12-21 Thursday\"),\n  @(\"43\u00d716=\", \"41\u00d755=\"),\n  @(\"17\u00d724=\", \"45\u00d719=\"),\n  @(\"41\u00d752=\", \"63\u00d762=\"),\n  @(\"26\u00d768=\", \"64\u00d724=\"),\n  @(\"83\u00d784=\", \"35\u00d793=\"),\n  @(\"89\u00d747=\", \"61\u00d724=\"),\n  @(\"46\u00d718=\", \"32\u00d735=\"),\n  @(\"69\u00d798=\", \"66\u00d776=\"),\n  @(\"16\u00d719=\", \"41\u00d738=\"),\n  @(\"35\u00d769=\", \"29\u00d796=\"),\n  @(\"82\u00d723=\", \"16\u00d776=\"),\n  @(\"58\u00d789=\", \"46\u00d778=\"),\n  @(\"70\u00d728=\", \"64\u00d714=\"),\n  @(\"67\u00d755=\", \"54\u00d775=\"),\n  @(\"27\u00d730=\", \"77\u00d784=\"),\n  @(\"94\u00d736=\", \"92\u00d797=\"),\n  @(\"89\u00d711=\", \"98\u00d718=\"),\n  @(\"88\u00d773=\", \"55\u00d712=\"),\n  @(\"93\u00d739=\", \"39\u00d776=\"),\n  @(\"73\u00d723=\", \"66\u00d729=\"),\n  @(\"72\u00d724=\", \"55\u00d750=\"),\n  @(\"81\u00d779=\", \"93\u00d783=\"),\n  @(\"46\u00d717=\", \"94\u00d791=\"),\n  @(\"34\u00d728=\", \"30\u00d741=\"),\n  @(\"56\u00d795=\", \"41\u00d790=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n\n  $found = $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n\n  if (-not $found) {\n    throw \"No match found for '$oldText'\"\n  }\n}\n"}
